# Update vm_pu.xlsx bus voltage magnitude results for the "380 kV" case.
# Columns B,C,D,E,F,I,J,K,L,M,N contain per-bus voltage magnitudes (p.u.)
# for rows 2-25 (time steps 0-23). Column A (index), column G (slack,
# always 1) and column H (unused) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colLetters = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @"
1.02,1.010484309763178,1.013129844277631,1.012768256453185,1.02299249721028,1.025564463961677,1.015738288697893,1.015991922530632,1.015631421819044,1.025825246014936,1.009325151751361
1.02,1.011614037404863,1.01417050772705,1.01373179448759,1.024267018433759,1.025543591503628,1.016499889652448,1.016836590618084,1.016399100089886,1.026905275428904,1.009584421370086
1.02,1.012345284487156,1.014844454671031,1.01435581510319,1.025090401850141,1.025527613580591,1.016992457397564,1.01738312379491,1.016895776884773,1.027602272531213,1.009751899987088
1.02,1.012652758877663,1.015127918463131,1.014618284600377,1.025436238050663,1.025520303761093,1.017199476115593,1.017612881167826,1.017104565206043,1.027894846844519,1.009822239520789
1.02,1.012704388623659,1.01517552129963,1.014662362086786,1.025494287078344,1.025519041627896,1.017234232153914,1.017651458151765,1.017139620837127,1.027943945383178,1.009834045823369
1.02,1.012349392745159,1.014848241791156,1.014359321717504,1.025095024163671,1.025527518236482,1.016995223814292,1.017386193843941,1.016898566781524,1.027606183664622,1.009752840136414
1.02,1.010866058001038,1.013481424029326,1.01309377552401,1.023423501003549,1.025557921341,1.015995725067276,1.016277387223835,1.015890875400787,1.026190631304045,1.009412832404199
1.02,1.008253999695243,1.011077244982634,1.010867899226553,1.020467925226902,1.025592601075909,1.014232629301508,1.014323326574419,1.014114701408479,1.023682017282031,1.008811502024735
1.02,1.006513737877013,1.00947733490595,1.009386770785786,1.018490647883195,1.025603068279385,1.013055956041703,1.013020451634356,1.01293023022205,1.021999994357436,1.008409137104556
1.02,1.005760430650382,1.008785229824752,1.008746082072227,1.017632811224697,1.025604610085929,1.012546133159121,1.012456244658848,1.012417250865919,1.021269367774712,1.008234557019918
1.02,1.005480652862979,1.008528250495199,1.008508198815126,1.017313920930066,1.025604734273553,1.012356714077093,1.012246664252412,1.012226692539931,1.020997633840438,1.008169656974818
1.02,1.005540664633325,1.008583368952887,1.00855922119807,1.017382335393541,1.025604727918276,1.012397347304749,1.012291620372338,1.012267568634073,1.021055937399837,1.00818358066563
1.02,1.005737303445942,1.008763985812034,1.008726416591226,1.017606456795743,1.025604629496225,1.012530476705577,1.012438920861365,1.012401499563776,1.021246913242388,1.008229193448556
1.02,1.005858463572088,1.008875282955289,1.008829444035713,1.017744511999903,1.025604509446884,1.012612495718652,1.012529676363024,1.012484016830383,1.021364533785618,1.008257289932278
1.02,1.006563737262973,1.009523281613565,1.009429304826826,1.018547544497377,1.025602903014827,1.013089784594274,1.013057894999519,1.012964272908698,1.022048435085864,1.008420715952515
1.02,1.007006199276307,1.009929931921766,1.009805755203822,1.019050819208078,1.025601095246419,1.013389090274921,1.013389217719945,1.013265498862141,1.02247681180164,1.008523134065088
1.02,1.007264302920423,1.010167188533625,1.010025395060942,1.019344210544008,1.025599752207639,1.01356363981453,1.013582467544463,1.013441189689116,1.022726454900511,1.008582838726477
1.02,1.007352313550602,1.010248097849866,1.01010029723146,1.019444222274182,1.025599245290644,1.013623151563128,1.013648359896832,1.013501094182843,1.02281153916261,1.008603190694662
1.02,1.00695872491843,1.009886295540187,1.009765359157303,1.018996839186369,1.025601319048141,1.013356980757649,1.01335367046189,1.013233181121353,1.022430873983433,1.008512149097958
1.02,1.005679397289765,1.008710795913887,1.008677179040332,1.017540465587766,1.025604670853805,1.012491274750777,1.012395544788256,1.012362060660743,1.021190685196861,1.008215763091283
1.02,1.004875229465402,1.007972286829223,1.007993558282928,1.01662332932573,1.025604183719274,1.011946691991626,1.011793081766879,1.011814265839881,1.020408922826244,1.008029105415367
1.02,1.005301515993101,1.008363730318997,1.008355905563232,1.017109659262828,1.025604687631517,1.012235412316294,1.01211246394293,1.012104670670096,1.020823540448983,1.008128085423114
1.02,1.006980176465585,1.009906012743518,1.009783612197847,1.019021230929136,1.025601218813544,1.013371489760013,1.013369732753195,1.013247784145415,1.022451631990768,1.00851711283702
1.02,1.008929077969227,1.011698273461991,1.011442848539452,1.021233221707093,1.025585869318646,1.01468865442251,1.014828525299383,1.014573945814485,1.024332245126138,1.008967220596926
"@

$rows = $data -split "`n" | Where-Object { $_.Trim().Length -gt 0 }

$rowIndex = 2
foreach ($line in $rows) {
    $values = $line -split ","
    for ($i = 0; $i -lt $colLetters.Length; $i++) {
        $cellRef = "$($colLetters[$i])$rowIndex"
        $ws.Range($cellRef).Value = [double]$values[$i]
    }
    $rowIndex++
}

Write-Output "Updated $($rows.Count) rows of vm_pu results"
